$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'XGBoost (Baseline)'
$ws.Range('B2').Value = '[''temp'', ''rhum'', ''dayofweek'', ''holiday'', ''timesofday'', ''wdsp'', ''rainfall_intensity'', ''peak'', ''working_day'', ''hour'', ''season'']'
$ws.Range('C2').Value = '{''objective'': ''reg:squarederror'', ''base_score'': 0.5, ''booster'': ''gbtree'', ''colsample_bylevel'': 1, ''colsample_bynode'': 1, ''colsample_bytree'': 1, ''gamma'': 0, ''gpu_id'': -1, ''importance_type'': ''gain'', ''interaction_constraints'': '''', ''learning_rate'': 0.300000012, ''max_delta_step'': 0, ''max_depth'': 3, ''min_child_weight'': 1, ''missing'': nan, ''monotone_constraints'': ''()'', ''n_estimators'': 500, ''n_jobs'': 12, ''num_parallel_tree'': 1, ''random_state'': 42, ''reg_alpha'': 0, ''reg_lambda'': 1, ''scale_pos_weight'': 1, ''subsample'': 1, ''tree_method'': ''exact'', ''validate_parameters'': 1, ''verbosity'': None, ''seed'': 42, ''eval_metric'': ''rmse''}'
$ws.Range('D2').Value = '[{ ''metric'': RSME, ''train'': 2.4729,  ''validation'': 2.6320, ''test'': None }, { ''metric'': MAE, ''train'': 1.8143,  ''validation'': 1.9269, ''test'': None }]'
$ws.Range('E2').Value = '01/07/2022 22:14:20'
$ws.Range('F2').Value = 'Baseline XGBoost'

$ws.Range('A3').Value = 'XGBoost 1'
$ws.Range('B3').Value = '[''temp'', ''rhum'', ''dayofweek'', ''holiday'', ''timesofday'', ''wdsp'', ''rainfall_intensity'', ''peak'', ''working_day'', ''hour'', ''season'']'
$ws.Range('C3').Value = '{''objective'': ''reg:squarederror'', ''base_score'': 0.5, ''booster'': ''gbtree'', ''colsample_bylevel'': 1, ''colsample_bynode'': 1, ''colsample_bytree'': 0.5, ''gamma'': 1, ''gpu_id'': -1, ''importance_type'': ''gain'', ''interaction_constraints'': '''', ''learning_rate'': 0.200000003, ''max_delta_step'': 0, ''max_depth'': 3, ''min_child_weight'': 1, ''missing'': nan, ''monotone_constraints'': ''()'', ''n_estimators'': 500, ''n_jobs'': 12, ''num_parallel_tree'': 1, ''random_state'': 42, ''reg_alpha'': 0, ''reg_lambda'': 1, ''scale_pos_weight'': 1, ''subsample'': 1, ''tree_method'': ''exact'', ''validate_parameters'': 1, ''verbosity'': None, ''eta'': 0.2, ''seed'': 42, ''eval_metric'': ''rmse''}'
$ws.Range('D3').Value = '[{ ''metric'': RSME, ''train'': 2.4930,  ''validation'': 2.6233, ''test'': None }, { ''metric'': MAE, ''train'': 1.8292,  ''validation'': 1.9204, ''test'': None }]'
$ws.Range('E3').Value = '01/07/2022 22:14:21'
$ws.Range('F3').Value = ''

$ws.Range('A4').Value = 'XGBoost 2'
$ws.Range('B4').Value = '[''temp'', ''rhum'', ''dayofweek'', ''holiday'', ''timesofday'', ''wdsp'', ''rainfall_intensity'', ''peak'', ''working_day'', ''hour'', ''season'']'
$ws.Range('C4').Value = '{''objective'': ''reg:squarederror'', ''base_score'': 0.5, ''booster'': ''gbtree'', ''colsample_bylevel'': 1, ''colsample_bynode'': 1, ''colsample_bytree'': 0.5, ''gamma'': 1, ''gpu_id'': -1, ''importance_type'': ''gain'', ''interaction_constraints'': '''', ''learning_rate'': 0.00999999978, ''max_delta_step'': 0, ''max_depth'': 9, ''min_child_weight'': 1, ''missing'': nan, ''monotone_constraints'': ''()'', ''n_estimators'': 1000, ''n_jobs'': 12, ''num_parallel_tree'': 1, ''random_state'': 42, ''reg_alpha'': 0, ''reg_lambda'': 1, ''scale_pos_weight'': 1, ''subsample'': 0.7, ''tree_method'': ''exact'', ''validate_parameters'': 1, ''verbosity'': None, ''eta'': 0.01, ''seed'': 42, ''eval_metric'': ''rmse''}'
$ws.Range('D4').Value = '[{ ''metric'': RSME, ''train'': 1.8872,  ''validation'': 2.5822, ''test'': None }, { ''metric'': MAE, ''train'': 1.3823,  ''validation'': 1.8758, ''test'': None }]'
$ws.Range('E4').Value = '01/07/2022 22:14:54'
$ws.Range('F4').Value = ''

$ws.Range('A5').Value = 'XGBoost (gamma: 1.5)'
$ws.Range('B5').Value = '[''temp'', ''rhum'', ''dayofweek'', ''holiday'', ''timesofday'', ''wdsp'', ''rainfall_intensity'', ''peak'', ''working_day'', ''hour'', ''season'']'
$ws.Range('C5').Value = '{''objective'': ''reg:squarederror'', ''base_score'': 0.5, ''booster'': ''gbtree'', ''colsample_bylevel'': 1, ''colsample_bynode'': 1, ''colsample_bytree'': 0.5, ''gamma'': 1.5, ''gpu_id'': -1, ''importance_type'': ''gain'', ''interaction_constraints'': '''', ''learning_rate'': 0.00999999978, ''max_delta_step'': 0, ''max_depth'': 9, ''min_child_weight'': 1, ''missing'': nan, ''monotone_constraints'': ''()'', ''n_estimators'': 1000, ''n_jobs'': 12, ''num_parallel_tree'': 1, ''random_state'': 42, ''reg_alpha'': 0, ''reg_lambda'': 1, ''scale_pos_weight'': 1, ''subsample'': 0.7, ''tree_method'': ''exact'', ''validate_parameters'': 1, ''verbosity'': None, ''eta'': 0.01, ''seed'': 42, ''eval_metric'': ''rmse''}'
$ws.Range('D5').Value = '[{ ''metric'': RSME, ''train'': 1.8869,  ''validation'': 2.5816, ''test'': None }, { ''metric'': MAE, ''train'': 1.3840,  ''validation'': 1.8753, ''test'': None }]'
$ws.Range('E5').Value = '01/07/2022 22:15:29'
$ws.Range('F5').Value = ''

$ws.Range('A6').Value = 'XGBoost (eta: 0.001)'
$ws.Range('B6').Value = '[''temp'', ''rhum'', ''dayofweek'', ''holiday'', ''timesofday'', ''wdsp'', ''rainfall_intensity'', ''peak'', ''working_day'', ''hour'', ''season'']'
$ws.Range('C6').Value = '{''objective'': ''reg:squarederror'', ''base_score'': 0.5, ''booster'': ''gbtree'', ''colsample_bylevel'': 1, ''colsample_bynode'': 1, ''colsample_bytree'': 0.5, ''gamma'': 1.5, ''gpu_id'': -1, ''importance_type'': ''gain'', ''interaction_constraints'': '''', ''learning_rate'': 0.00100000005, ''max_delta_step'': 0, ''max_depth'': 11, ''min_child_weight'': 1, ''missing'': nan, ''monotone_constraints'': ''()'', ''n_estimators'': 5000, ''n_jobs'': 12, ''num_parallel_tree'': 1, ''random_state'': 42, ''reg_alpha'': 0, ''reg_lambda'': 1, ''scale_pos_weight'': 1, ''subsample'': 0.7, ''tree_method'': ''exact'', ''validate_parameters'': 1, ''verbosity'': None, ''eta'': 0.001, ''seed'': 42, ''eval_metric'': ''rmse''}'
$ws.Range('D6').Value = '[{ ''metric'': RSME, ''train'': 1.6972,  ''validation'': 2.5873, ''test'': None }, { ''metric'': MAE, ''train'': 1.2327,  ''validation'': 1.8745, ''test'': None }]'
$ws.Range('E6').Value = '01/07/2022 22:21:59'
$ws.Range('F6').Value = ''

$ws.Range('A7').Value = 'XGBoost (max_depth: 9 and subsample: 0.8)'
$ws.Range('B7').Value = '[''temp'', ''rhum'', ''dayofweek'', ''holiday'', ''timesofday'', ''wdsp'', ''rainfall_intensity'', ''peak'', ''working_day'', ''hour'', ''season'']'
$ws.Range('C7').Value = '{''objective'': ''reg:squarederror'', ''base_score'': 0.5, ''booster'': ''gbtree'', ''colsample_bylevel'': 1, ''colsample_bynode'': 1, ''colsample_bytree'': 0.5, ''gamma'': 1.5, ''gpu_id'': -1, ''importance_type'': ''gain'', ''interaction_constraints'': '''', ''learning_rate'': 0.00100000005, ''max_delta_step'': 0, ''max_depth'': 9, ''min_child_weight'': 1, ''missing'': nan, ''monotone_constraints'': ''()'', ''n_estimators'': 5000, ''n_jobs'': 12, ''num_parallel_tree'': 1, ''random_state'': 42, ''reg_alpha'': 0, ''reg_lambda'': 1, ''scale_pos_weight'': 1, ''subsample'': 0.8, ''tree_method'': ''exact'', ''validate_parameters'': 1, ''verbosity'': None, ''eta'': 0.001, ''seed'': 42, ''eval_metric'': ''rmse''}'
$ws.Range('D7').Value = '[{ ''metric'': RSME, ''train'': 1.9097,  ''validation'': 2.5784, ''test'': None }, { ''metric'': MAE, ''train'': 1.3968,  ''validation'': 1.8723, ''test'': None }]'
$ws.Range('E7').Value = '01/07/2022 22:26:55'
$ws.Range('F7').Value = ''

